$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$range = $ws.Range("A2:A74")
$range.Style = "Normal"
